# A new daily price record was reported for "Choclo" (Dulce o Americano,
# Femacal de La Calera) and needs to be inserted into the historical log.
# The record belongs right after the existing row for 2023-07-08
# (serial 44901, row 1046), so insert a new row at 1047 - this pushes every
# following record down by one (old row 1047 becomes 1048, ..., old row
# 1131 becomes the new row 1132) - and populate it with the new reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift rows 1047..1131 down to 1048..1132, leaving a blank row 1047.
$ws.Rows.Item(1047).Insert()

# Fill the newly-opened row 1047 with the new observation.
$ws.Range("A1047").Value = 3
$ws.Range("B1047").Value = "Femacal de La Calera"
$ws.Range("C1047").Value = "Coquimbo"
$ws.Range("D1047").Value = 45132
$ws.Range("E1047").Value = 5
$ws.Range("F1047").Value = 100112024
$ws.Range("G1047").Value = "Choclo"
$ws.Range("H1047").Value = "Dulce o Americano"
$ws.Range("I1047").Value = "Primera"
$ws.Range("J1047").Value = 70
$ws.Range("K1047").Value = 32000
$ws.Range("L1047").Value = 33000
$ws.Range("M1047").Value = 32500
$ws.Range("N1047").Value = "$/malla 70 unidades"
$ws.Range("O1047").Value = "Región de Arica y Parinacota"
$ws.Range("P1047").Value = 464
$ws.Range("Q1047").Value = 70
$ws.Range("R1047").Value = "Hortaliza"
